# Melb House Price Prediction with Machine Learning - Group 7
# "moving html file to static and fixing presentation"
#
# Slide 7 ("Webpage") shuffles several picture/placeholder positions and
# replaces one image (previously shown twice, as "Picture 23" and
# "Picture 27") with three rotated copies of the "Picture 27" image.
#
# NOTE: PowerPoint's Shape.Left/Top/Width/Height/Rotation COM properties
# are 32-bit (Single) floats. The literal point values below were chosen
# so that, after the Single-precision truncation PowerPoint performs
# internally (points -> EMU, 1 pt = 12700 EMU), they land exactly on the
# target EMU offsets/extents from the authoritative OOXML.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# --- Content Placeholder 12 (the screenshot shown top-ish) ---------------
$shContent = $s.Shapes.Item("Content Placeholder 12")
$shContent.Left = 267.2547302246094
$shContent.Top  = 261.4000244140625

# --- Picture 4 (small inset image, upper area) ----------------------------
$shPic4 = $s.Shapes.Item("Picture 4")
$shPic4.Left = 90.15559387207031
$shPic4.Top  = 154.628662109375

# --- Picture 6 (small inset image, lower area) -----------------------------
$shPic6 = $s.Shapes.Item("Picture 6")
$shPic6.Left = 85.29551696777344
$shPic6.Top  = 360.1234130859375

# --- Picture 14 (right-hand side image) ------------------------------------
$shPic14 = $s.Shapes.Item("Picture 14")
$shPic14.Left = 487.158203125
$shPic14.Top  = 254.32142639160156

# --- Picture 23 is removed outright (its duplicate, Picture 27, is kept
#     and multiplied into three rotated copies instead). --------------------
$shPic23 = $s.Shapes.Item("Picture 23")
$shPic23.Delete()

# --- Picture 27 becomes a rotated, slightly resized copy --------------------
$shPic27 = $s.Shapes.Item("Picture 27")
$shPic27.Width    = 77.55464935302734
$shPic27.Height   = 113.04000091552734
$shPic27.Left     = 203.66213989257812
$shPic27.Top      = 181.88299560546875
$shPic27.Rotation = 39.3961

# --- Two more rotated copies of the same picture, appended after it --------
$shNew2 = $shPic27.Duplicate()
$shNew2.Name     = "Picture 2"
$shNew2.Width    = 75.61409759521484
$shNew2.Height   = 113.04000091552734
$shNew2.Left     = 219.757568359375
$shNew2.Top      = 332.35040283203125
$shNew2.Rotation = 337.4685666666667

$shNew3 = $shPic27.Duplicate()
$shNew3.Name     = "Picture 3"
$shNew3.Width    = 75.61409759521484
$shNew3.Height   = 113.04000091552734
$shNew3.Left     = 395.532470703125
$shNew3.Top      = 258.4296875
$shNew3.Rotation = 180.0
